$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update commit-driven data values (Public Ads Paging scores)
$ws.Range("C9").Value = 33
$ws.Range("C14").Value = 5

# Update view: scroll position (topLeftCell) and active selection
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("C15").Select()
